# Update EUR->ARS rate: append the latest quote as a new row at the
# bottom of the sheet (mirrors the existing Fecha / Hora / Cotizacion
# columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row after the existing data (row 95 -> 96).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Column A holds a date-looking string ("2025-10-23"). A bare .Value
# assignment would be auto-coerced into a date serial by Excel, so enter
# it with a leading apostrophe to force literal text, exactly like the
# rest of the column.
$ws.Cells.Item($newRow, 1).Value = "'2025-10-23"
$ws.Cells.Item($newRow, 2).Value = "21:20:30"
$ws.Cells.Item($newRow, 3).Value = "1.00 EUR = 1,826.0714"
